$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates (Q, R) to whole numbers for rows 2 and 3.
$ws.Range("Q2").Value = [Math]::Round($ws.Range("Q2").Value2, 0)
$ws.Range("R2").Value = [Math]::Round($ws.Range("R2").Value2, 0)
$ws.Range("Q3").Value = [Math]::Round($ws.Range("Q3").Value2, 0)
$ws.Range("R3").Value = [Math]::Round($ws.Range("R3").Value2, 0)

# Clear the now-unused start/end time columns for rows 2 and 3.
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
